$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-10 Saturday" "2024-08-11 Sunday"

Replace-Text "753×3=" "288×8="
Replace-Text "347×4=" "449×9="
Replace-Text "723×8=" "675×5="
Replace-Text "588×8=" "433×6="
Replace-Text "833×9=" "541×4="

Replace-Text "160×2=" "250×8="
Replace-Text "250×9=" "883×2="
Replace-Text "866×9=" "211×2="
Replace-Text "630×8=" "746×6="
Replace-Text "881×4=" "306×2="

Replace-Text "512×9=" "385×6="
Replace-Text "108×6=" "163×8="
Replace-Text "709×8=" "755×8="
Replace-Text "427×2=" "538×7="
Replace-Text "221×6=" "255×7="

Replace-Text "299×8=" "625×8="
Replace-Text "701×3=" "587×9="
Replace-Text "260×2=" "817×6="
Replace-Text "546×4=" "292×7="
Replace-Text "508×8=" "679×5="

Replace-Text "212×6=" "982×4="
Replace-Text "617×3=" "692×8="
Replace-Text "538×3=" "881×3="
Replace-Text "740×2=" "638×9="
Replace-Text "805×7=" "504×2="
